# schedule.xlsx: row 8 ("6. hét" week) gets the same "csúszás" (slip/delay)
# marker in column D that rows 6 and 7 already carry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D8").Value = "csúszás"
